$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the formatting of row 4 onto the still-empty row 5 so the new
#    row inherits the same borders/fonts/alignment as the other data rows.
$ws.Range("A4:M4").Copy()
$ws.Range("A5:M5").PasteSpecial(-4122)

# 2. Clear the old precondition text in F4 first -- it is only referenced
#    by this one cell, so clearing it drops the now-unused shared string
#    before we give F4 its replacement text.
$ws.Range("F4").Value = ""

# 3. Populate the new row (CP_AUTO_004 - "Redirigir ONT") and update F4,
#    in the exact order the new case's strings should be introduced.
$ws.Range("B5").Value = "Redirigir ONT"
$ws.Range("E5").Value = "Validar abrir modal redirigir ont y dar clic en la opcion ""NO"""
$ws.Range("F4").Value = "El usuario debe haber seleccionado la opcion configuracion de wifi en la lista de opciones"
$ws.Range("F5").Value = "El usuario debe haber seleccionado la opcion Redirigir ONTen la lista de opciones"
$ws.Range("A5").Value = "CP_AUTO_004"
$ws.Range("G5").Value = "1. Clic en el boton OPCIONES`n2. Clic en opción " + [char]8220 + "Redirigir ONT" + [char]8221 + "`n3.Clic en el botón ""NO"" del modal de confirmación"
$ws.Range("I5").Value = "El sistema debe permitir cerrar el modal una vez se selecciona ""NO"""
$ws.Range("J5").Value = "El modal de redirigir ont se cierra correctamente."

# 4. Remaining row 5 cells reuse existing strings.
$ws.Range("C5").Value = "Positivo"
$ws.Range("D5").Value = "eCenter"
$ws.Range("H5").Value = "N/A"
$ws.Range("K5").Value = "OK"
$ws.Range("L5").Value = "SI"
$ws.Range("M5").Value = "N/A"

# 5. Row heights: row 4 grew to fit its new precondition text, row 5 is
#    the newly added case.
$ws.Rows.Item(4).RowHeight = 182.25
$ws.Rows.Item(5).RowHeight = 144

# 6. Match the saved selection/active cell.
$ws.Range("M5").Select()
